$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the footnote rows (A8:A13) down to (A18:A23), leaving the blank
# styled row (A16) where it is, and opening a gap above it.
$src = $ws.Range("A8:A13")
$dst = $ws.Range("A18")
$src.Cut($dst)

# Match the author's last selection on the sheet.
$ws.Range("B9").Select()
